$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Hiaral50"
$ws.Range("C5").Value = "'$1.40"
$ws.Range("D5").Value = "'$4.19"

$ws.Range("A6").Value = "Raesty92"
$ws.Range("C6").Value = "'$2.91"
$ws.Range("D6").Value = "'$8.73"
